# Update betting-odds values on Sheet1 for the 2024-11-17 FlashScore export.
# Each line below sets exactly one cell to its new value, matching the
# canonical-XML diff (rows 2,3,6,7,8,14,15,17,18,20).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("L2").Value = 3.75
$ws.Range("N2").Value = 4.75
$ws.Range("AH2").Value = 11
$ws.Range("AQ2").Value = 81
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44
$ws.Range("G6").Value = 1.8
$ws.Range("H6").Value = 3.25
$ws.Range("K6").Value = 2.07
$ws.Range("L6").Value = 4.8
$ws.Range("O6").Value = 1.35
$ws.Range("P6").Value = 2.72
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.65
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.5
$ws.Range("U6").Value = 1.87
$ws.Range("V6").Value = 1.75
$ws.Range("W6").Value = 6.3
$ws.Range("AA6").Value = 15
$ws.Range("AB6").Value = 30
$ws.Range("AC6").Value = 8.25
$ws.Range("AD6").Value = 6.4
$ws.Range("AE6").Value = 16.5
$ws.Range("AF6").Value = 90
$ws.Range("AG6").Value = 10.5
$ws.Range("AI6").Value = 15
$ws.Range("AJ6").Value = 80
$ws.Range("AK6").Value = 50
$ws.Range("AL6").Value = 55
$ws.Range("AM6").Value = 800
$ws.Range("AN6").Value = 3.55
$ws.Range("AO6").Value = 8.75
$ws.Range("AP6").Value = 17.5
$ws.Range("AQ6").Value = 30
$ws.Range("AS6").Value = 250
$ws.Range("AT6").Value = 2.47
$ws.Range("AU6").Value = 7.3
$ws.Range("AV6").Value = 70
$ws.Range("AX6").Value = 26
$ws.Range("AY6").Value = 32
$ws.Range("AZ6").Value = 175
$ws.Range("BA6").Value = 200
$ws.Range("BB6").Value = 450
$ws.Range("G7").Value = 2
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 2.6
$ws.Range("AE7").Value = 13
$ws.Range("AJ7").Value = 41
$ws.Range("AK7").Value = 26
$ws.Range("AM7").Value = 151
$ws.Range("AX7").Value = 19
$ws.Range("BB7").Value = 151
$ws.Range("G8").Value = 2.38
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 2.25
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 1.73
$ws.Range("R8").Value = 2.08
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.2
$ws.Range("AM8").Value = 151
$ws.Range("AO8").Value = 13
$ws.Range("Q14").Value = 1.6
$ws.Range("R14").Value = 2.3
$ws.Range("G15").Value = 3
$ws.Range("I15").Value = 2.63
$ws.Range("J15").Value = 3.75
$ws.Range("L15").Value = 3.5
$ws.Range("Y15").Value = 12
$ws.Range("AA15").Value = 29
$ws.Range("AG15").Value = 6.5
$ws.Range("AR15").Value = 101
$ws.Range("AW15").Value = 4.5
$ws.Range("AX15").Value = 17
$ws.Range("M17").Value = 1.08
$ws.Range("N17").Value = 8
$ws.Range("O18").Value = 1.53
$ws.Range("P18").Value = 2.38
$ws.Range("U18").Value = 2.2
$ws.Range("V18").Value = 1.62
$ws.Range("AE18").Value = 19
$ws.Range("AG18").Value = 7
$ws.Range("BA18").Value = 101
$ws.Range("M20").Value = 1.04
$ws.Range("N20").Value = 6.3
